# Update attendance/sales figures (column F) on multiple sheets to match
# data output generated at commit 456a3b4 (gh-pages update).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7958
$ws.Range("F7").Value = 5825
$ws.Range("F9").Value = 2849
$ws.Range("F10").Value = 1235
$ws.Range("F14").Value = 658
$ws.Range("F16").Value = 4096
$ws.Range("F17").Value = 4096
$ws.Range("F20").Value = 81
$ws.Range("F23").Value = 40
$ws.Range("F24").Value = 6030
$ws.Range("F25").Value = 6030
$ws.Range("F29").Value = 411
$ws.Range("F31").Value = 434
$ws.Range("F32").Value = 4304
$ws.Range("F33").Value = 1580
$ws.Range("F37").Value = 89
$ws.Range("F39").Value = 68
$ws.Range("F40").Value = 46
$ws.Range("F41").Value = 3800
$ws.Range("F42").Value = 50
$ws.Range("F43").Value = 61
$ws.Range("F45").Value = 2371
$ws.Range("F50").Value = 246

# --- Sheet "演出" (Performances) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 35

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1389

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1389
$ws.Range("F4").Value = 7958
$ws.Range("F7").Value = 5826
$ws.Range("F9").Value = 2849
$ws.Range("F10").Value = 1235
$ws.Range("F15").Value = 658
$ws.Range("F17").Value = 4096
$ws.Range("F18").Value = 4096
$ws.Range("F21").Value = 81
$ws.Range("F24").Value = 40
$ws.Range("F25").Value = 6030
$ws.Range("F26").Value = 6031
$ws.Range("F29").Value = 411
$ws.Range("F30").Value = 184
$ws.Range("F31").Value = 434
$ws.Range("F33").Value = 4304
$ws.Range("F34").Value = 1580
$ws.Range("F40").Value = 89
$ws.Range("F42").Value = 3800
$ws.Range("F43").Value = 61
$ws.Range("F46").Value = 35
$ws.Range("F47").Value = 2371
$ws.Range("F50").Value = 246
